$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "633÷7=90, 3"; New = "305÷9=33, 8" },
    @{ Old = "250÷6=41, 4"; New = "698÷2=349, 0" },
    @{ Old = "407÷4=101, 3"; New = "498÷3=166, 0" },
    @{ Old = "523÷4=130, 3"; New = "654÷4=163, 2" },
    @{ Old = "885÷2=442, 1"; New = "878÷6=146, 2" },
    @{ Old = "762÷5=152, 2"; New = "883÷4=220, 3" },
    @{ Old = "972÷4=243, 0"; New = "134÷9=14, 8" },
    @{ Old = "528÷6=88, 0"; New = "117÷3=39, 0" },
    @{ Old = "493÷7=70, 3"; New = "934÷4=233, 2" },
    @{ Old = "134÷8=16, 6"; New = "655÷5=131, 0" },
    @{ Old = "475÷8=59, 3"; New = "167÷3=55, 2" },
    @{ Old = "625÷9=69, 4"; New = "484÷8=60, 4" },
    @{ Old = "458÷6=76, 2"; New = "490÷4=122, 2" },
    @{ Old = "984÷9=109, 3"; New = "302÷4=75, 2" },
    @{ Old = "568÷5=113, 3"; New = "918÷3=306, 0" },
    @{ Old = "345÷2=172, 1"; New = "870÷9=96, 6" },
    @{ Old = "387÷9=43, 0"; New = "500÷2=250, 0" },
    @{ Old = "900÷7=128, 4"; New = "791÷7=113, 0" },
    @{ Old = "858÷4=214, 2"; New = "957÷4=239, 1" },
    @{ Old = "512÷6=85, 2"; New = "533÷4=133, 1" },
    @{ Old = "268÷2=134, 0"; New = "332÷6=55, 2" },
    @{ Old = "296÷2=148, 0"; New = "288÷2=144, 0" },
    @{ Old = "623÷3=207, 2"; New = "599÷7=85, 4" },
    @{ Old = "389÷3=129, 2"; New = "170÷8=21, 2" },
    @{ Old = "713÷6=118, 5"; New = "159÷3=53, 0" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
